$wb = $excel.ActiveWorkbook

# --- Sheet "AMSIN": add a new row 29 after the existing data ---
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$wsAmsin.Range("A29").Value = "2023-07-31"
$wsAmsin.Range("B29").Value = 45138.43014296296
$wsAmsin.Range("C29").Value = "180aadhar"
$wsAmsin.Range("D29").Value = 33
$wsAmsin.Range("E29").Value = 33
$wsAmsin.Range("F29").Value = 0
$wsAmsin.Range("G29").Value = 1.59

# Copy the style down from the previous row (A28:G28) so the new row29
# (and the now-styled row28) match the expected formatting.
$wsAmsin.Range("A28:G28").Copy() | Out-Null
$wsAmsin.Range("A28:G29").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Re-set the values since PasteSpecial(formats) shouldn't touch values, but
# make sure nothing was altered.
$wsAmsin.Range("A28").Value = "2023-06-12"
$wsAmsin.Range("B28").Value = 45089.6111671875
$wsAmsin.Range("C28").Value = "178aadhdy"
$wsAmsin.Range("D28").Value = 33
$wsAmsin.Range("E28").Value = 33
$wsAmsin.Range("F28").Value = 0
$wsAmsin.Range("G28").Value = 1.47

$wsAmsin.Range("A29").Value = "2023-07-31"
$wsAmsin.Range("B29").Value = 45138.43014296296
$wsAmsin.Range("C29").Value = "180aadhar"
$wsAmsin.Range("D29").Value = 33
$wsAmsin.Range("E29").Value = 33
$wsAmsin.Range("F29").Value = 0
$wsAmsin.Range("G29").Value = 1.59

# --- Sheet "AMS": add new rows 17 and 18 after the existing data ---
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Range("A16:G16").Copy() | Out-Null
$wsAms.Range("A17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wsAms.Range("A17").Value = "2023-08-01"
$wsAms.Range("B17").Value = 45139.53672981481
$wsAms.Range("C17").Value = "180adhara"
$wsAms.Range("D17").Value = 33
$wsAms.Range("E17").Value = 33
$wsAms.Range("F17").Value = 0
$wsAms.Range("G17").Value = 1.5

$wsAms.Range("A18").Value = "2023-08-01"
$wsAms.Range("B18").Value = 45139.8615866535
$wsAms.Range("C18").Value = "180liveaadhar"
$wsAms.Range("D18").Value = 33
$wsAms.Range("E18").Value = 33
$wsAms.Range("F18").Value = 0
$wsAms.Range("G18").Value = 1.11
